$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Results")

# New rows of data to append after the existing last row (186)
$newRows = @(
    @(800, 100, 0.5, 0.6, 3397, 0),
    @(800, 100, 0.5, 0.6, 5887, 0),
    @(800, 100, 0.5, 0.6, 11178, -100)
)

$startRow = 187
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
    $ws.Cells.Item($r, 4).Value = $rowData[3]
    $ws.Cells.Item($r, 5).Value = $rowData[4]
    $ws.Cells.Item($r, 6).Value = $rowData[5]
}
